$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update category labels from 액티비티 to 관광지
$ws.Range("C8").Value = "관광지"
$ws.Range("C10").Value = "관광지"
$ws.Range("C12").Value = "관광지"

# Update category labels from 포인트 to 관광지
$ws.Range("C21").Value = "관광지"
$ws.Range("C22").Value = "관광지"
$ws.Range("C23").Value = "관광지"
$ws.Range("C25").Value = "관광지"
$ws.Range("C26").Value = "관광지"

# Update the Solitaire map URL and coordinates in row 19
$ws.Range("E19").Value = "https://www.google.com/maps/place/%EC%86%94%EB%A6%AC%ED%85%8C%EC%96%B4/@-23.8931308,15.9949915,15.86z/data=!4m6!3m5!1s0x1c731a93fc9fae0d:0xc238afe9f76ca6f7!8m2!3d-23.8933454!4d16.0045091!16zL20vMGdidHR2?entry=ttu&g_ep=EgoyMDI0MTEyNC4xIKXMDSoASAFQAw%3D%3D"
$ws.Range("F19").Value = -23.8931308
$ws.Range("G19").Value = 15.9949915
